$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.299.48'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.183.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.21'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.55'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.181.59'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.27'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.01%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.45'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +9.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.701.45'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.344.95'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.45%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +7.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.171.44'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.88%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '515.28'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.06'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.18'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +13.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.743'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +9.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.92'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.20'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.67%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.25'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +16.72%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +9.21%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +14.99%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +8.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.34'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +12.68%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +7.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.99'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '480.72'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +8.01%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.17'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +11.83%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0880'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +8.65%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.133.16'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.71'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.122'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.55'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +18.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.294'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +12.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.52'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0593'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +15.38%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.35'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +12.94%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.116'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '125.18'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.16%  '
